$d = $word.ActiveDocument

# 0x00B050 in RGB -> BGR-ordered long value Word COM expects for Font.Color
$green = 5287936  # RGB(0x00,0xB0,0x50) == wdColor 0x0050B000

$targets = @(
    "Both of the following areas should have progress - 5% ",
    "Queries should be written in LINQ"
)

foreach ($t in $targets) {
    $rng = $d.Content
    $found = $rng.Find.Execute($t, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $para = $rng.Paragraphs(1)
        $para.Range.Font.Color = $green
    }
}
